# Updates cryptos price/volume data per the scheduled GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.665.20'
$ws.Range("E2").Value = '  +0.60%  '

# Row 3
$ws.Range("D3").Value = '1.849.25'
$ws.Range("E3").Value = '  +0.04%  '

# Row 4
$ws.Range("D4").Value = '''0.9994'
$ws.Range("E4").Value = '  +0.03%  '

# Row 5
$ws.Range("D5").Value = '''262.92'
$ws.Range("E5").Value = '  -1.27%  '

# Row 6
$ws.Range("E6").Value = '  +0.05%  '

# Row 7
$ws.Range("D7").Value = '''0.5326'
$ws.Range("E7").Value = '  +1.91%  '

# Row 8
$ws.Range("D8").Value = '''0.3206'
$ws.Range("E8").Value = '  -2.49%  '

# Row 9
$ws.Range("D9").Value = '''0.06974'
$ws.Range("E9").Value = '  +2.17%  '

# Row 10
$ws.Range("D10").Value = '''19.16'
$ws.Range("E10").Value = '  +1.15%  '

# Row 11
$ws.Range("D11").Value = '''0.7814'
$ws.Range("E11").Value = '  -0.01%  '

# Row 12
$ws.Range("D12").Value = '''0.07832'
$ws.Range("E12").Value = '  +0.87%  '

# Row 13
$ws.Range("D13").Value = '1.833.06'
$ws.Range("E13").Value = '  -0.92%  '

# Row 14
$ws.Range("D14").Value = '''89.46'
$ws.Range("E14").Value = '  +1.11%  '

# Row 15
$ws.Range("D15").Value = '''5.062'

# Row 16
$ws.Range("D16").Value = '''14.20'
$ws.Range("E16").Value = '  +1.51%  '

# Row 17
$ws.Range("D17").Value = '''0.9991'
$ws.Range("E17").Value = '  +0.08%  '

# Row 18
$ws.Range("D18").Value = '''0.000007994'
$ws.Range("E18").Value = '  -0.01%  '

# Row 20
$ws.Range("D20").Value = '26.680.82'
$ws.Range("E20").Value = '  +0.55%  '

# Row 21
$ws.Range("D21").Value = '2.068.25'
$ws.Range("E21").Value = '  -0.53%  '

# Row 22
$ws.Range("D22").Value = '''4.655'
$ws.Range("E22").Value = '  +0.18%  '

# Row 23
$ws.Range("D23").Value = '''6.041'
$ws.Range("E23").Value = '  +0.42%  '

# Row 24
$ws.Range("D24").Value = '''9.440'
$ws.Range("E24").Value = '  -1.57%  '

# Row 25
$ws.Range("D25").Value = '''2.228'
$ws.Range("E25").Value = '  +1.34%  '

# Row 26
$ws.Range("D26").Value = '''142.68'
$ws.Range("E26").Value = '  -1.21%  '

# Row 27
$ws.Range("D27").Value = '''1.701'
$ws.Range("E27").Value = '  +2.39%  '

# Row 28
$ws.Range("D28").Value = '''17.12'
$ws.Range("E28").Value = '  +0.53%  '

# Row 29
$ws.Range("D29").Value = '''111.80'
$ws.Range("E29").Value = '  -0.21%  '

# Row 30
$ws.Range("D30").Value = '''4.314'
$ws.Range("E30").Value = '  +2.84%  '

# Row 31
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").Value = '''4.133'
$ws.Range("E31").Value = '  -0.51%  '

# Row 32
$ws.Range("B32").Value = 'Stellar'
$ws.Range("C32").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D32").Value = '''0.08770'
$ws.Range("E32").Value = '  +0.25%  '

# Row 33
$ws.Range("D33").Value = '''0.04886'
$ws.Range("E33").Value = '  +0.83%  '

# Row 34
$ws.Range("D34").Value = '''0.7391'
$ws.Range("E34").Value = '  +1.63%  '

# Row 35
$ws.Range("D35").Value = '''1.146'
$ws.Range("E35").Value = '  +0.55%  '

# Row 36
$ws.Range("D36").Value = '''2.865'
$ws.Range("E36").Value = '  +0.63%  '

# Row 37
$ws.Range("D37").Value = '''3.113'
$ws.Range("E37").Value = '  +0.22%  '

# Row 38
$ws.Range("D38").Value = '''2.383'
$ws.Range("E38").Value = '  +6.50%  '

# Row 39
$ws.Range("D39").Value = '''0.01755'
$ws.Range("E39").Value = '  -1.90%  '

# Row 40
$ws.Range("D40").Value = '''0.4856'
$ws.Range("E40").Value = '  -1.38%  '

# Row 41
$ws.Range("D41").Value = '''0.9094'
$ws.Range("E41").Value = '  -0.70%  '

# Row 42
$ws.Range("D42").Value = '''109.68'
$ws.Range("E42").Value = '  -1.67%  '

# Row 43
$ws.Range("D43").Value = '''5.926'
$ws.Range("E43").Value = '  -2.92%  '

# Row 44
$ws.Range("E44").Value = '  +0.12%  '

# Row 45
$ws.Range("D45").Value = '''7.767'
$ws.Range("E45").Value = '  -0.40%  '

# Row 46
$ws.Range("D46").Value = '''0.4232'
$ws.Range("E46").Value = '  +0.56%  '

# Row 47
$ws.Range("D47").Value = '''0.1260'
$ws.Range("E47").Value = '  +0.77%  '

# Row 48
$ws.Range("D48").Value = '''9.101'
$ws.Range("E48").Value = '  -0.23%  '

# Row 49
$ws.Range("D49").Value = '''35.13'
$ws.Range("E49").Value = '  +0.04%  '

# Row 50
$ws.Range("D50").Value = '''0.05843'

# Row 51
$ws.Range("D51").Value = '''0.8982'
$ws.Range("E51").Value = '  +0.89%  '

"Cryptos list updated"